$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-05-01 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-02 Thursday", 2) | Out-Null

# Update the 100 equation cells (20 rows x 5 columns), addressed positionally
# to avoid any ambiguity from repeated/overlapping equation text.
$t = $d.Tables.Item(1)
$values = @(
    @("64+3=", "23+17=", "99-14=", "96+2=", "71+22="),
    @("40+2=", "3+65=", "74+14=", "68-5=", "97-1="),
    @("18+6=", "39-11=", "5+73=", "78-66=", "20-9="),
    @("21+12=", "81-60=", "95-40=", "58-34=", "92-83="),
    @("92-8=", "54-26=", "97-67=", "7+15=", "94-54="),
    @("77+2=", "75-51=", "69-49=", "53-4=", "66-10="),
    @("85-43=", "37+19=", "3+79=", "78-9=", "98-94="),
    @("50+48=", "76+9=", "52-35=", "75+22=", "20+16="),
    @("38-28=", "76+19=", "83-62=", "86-55=", "21+58="),
    @("1+8=", "87+5=", "70-5=", "56-41=", "68-40="),
    @("14+67=", "28+60=", "38+3=", "26+39=", "34+58="),
    @("89-65=", "68-14=", "0+81=", "46-9=", "15+61="),
    @("24+17=", "93-87=", "66+26=", "14-4=", "36+46="),
    @("1+8=", "42-22=", "22+46=", "20+9=", "8+72="),
    @("22+46=", "83-41=", "16+52=", "74-9=", "62+9="),
    @("54+20=", "32+66=", "31+35=", "63+33=", "35-33="),
    @("33-18=", "18+39=", "15-13=", "57+29=", "57+14="),
    @("17+74=", "99-0=", "58-24=", "67+26=", "12+82="),
    @("98-5=", "84+15=", "50+7=", "14+49=", "33-6="),
    @("24+6=", "1+69=", "74+17=", "56+17=", "38-17=")
)

for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$r-1][$c-1]
    }
}

Write-Host "Replacements applied"
